$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 780, pushing existing rows 780:886 down to 781:887
$ws.Rows.Item(780).Insert()

# Populate the new row 780 with the new price observation
$ws.Cells.Item(780, 1).Value = 10
$ws.Cells.Item(780, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(780, 3).Value = "La Araucanía"
$ws.Cells.Item(780, 4).Value = 45154
$ws.Cells.Item(780, 5).Value = 9
$ws.Cells.Item(780, 6).Value = 100112045
$ws.Cells.Item(780, 7).Value = "Zapallo"
$ws.Cells.Item(780, 8).Value = "Camote"
$ws.Cells.Item(780, 9).Value = "1a (guarda)"
$ws.Cells.Item(780, 10).Value = 100
$ws.Cells.Item(780, 11).Value = 500
$ws.Cells.Item(780, 12).Value = 500
$ws.Cells.Item(780, 13).Value = 500
$ws.Cells.Item(780, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(780, 15).Value = "Región del Maule"
$ws.Cells.Item(780, 16).Value = 500
$ws.Cells.Item(780, 17).Value = 1
$ws.Cells.Item(780, 18).Value = "Hortaliza"
